# Daily attendance processing - reorders the "Recorded By" (column G) list
# so that the "System" entry is surfaced first. For every data row whose
# G-cell holds a comma-separated list of recorders, the first (oldest)
# recorder is rotated to the end of the list, unless the list already
# starts with "System" (case-insensitive) - in that case it is left as-is.
# Single-value cells are untouched since there's nothing to reorder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow  = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $raw  = $cell.Value2

    if ($raw -eq $null) { continue }
    if ($raw.GetType().Name -ne "String") { continue }
    if ($raw.IndexOf(",") -lt 0) { continue }

    $parts = $raw.Split(",")
    $trimmedParts = @()
    foreach ($p in $parts) {
        $trimmedParts += $p.Trim()
    }

    if ($trimmedParts[0].ToLower() -eq "system") {
        continue
    }

    $rotated = $trimmedParts[1..($trimmedParts.Count - 1)] + $trimmedParts[0]
    $newValue = [string]::Join(", ", $rotated)

    $cell.Value2 = $newValue
}
